# Edit: 20250502_trend_summary.xlsx
# - "Summary Table": remove the placeholder "---------" row (old row 2), re-style header
# - "Cooccurrence": populate with source/target/count keyword co-occurrence data
# - "Associations": replace placeholder row and append the remaining association rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary Table
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary Table")

# Remove the placeholder dashed-line row; rows below shift up automatically,
# preserving their original (already-correct) values.
$wsSummary.Rows.Item(2).Delete()

# Re-set the header labels without the stray leading/trailing spaces.
$wsSummary.Cells.Item(1, 1).Value = "Keyword"
$wsSummary.Cells.Item(1, 2).Value = "Keyword Count"
$wsSummary.Cells.Item(1, 3).Value = "Short Summary"
$wsSummary.Cells.Item(1, 4).Value = "Source URL"
$wsSummary.Cells.Item(1, 5).Value = "Detailed Summary"

# ---------------------------------------------------------------------------
# 2) Cooccurrence
# ---------------------------------------------------------------------------
$wsCooccurrence = $wb.Worksheets.Item("Cooccurrence")

$wsCooccurrence.Cells.Item(1, 1).Value = "source"
$wsCooccurrence.Cells.Item(1, 2).Value = "target"
$wsCooccurrence.Cells.Item(1, 3).Value = "count"

$cooccurrenceData = @(
    @("人工智能", "新质生产力", 1),
    @("创新驱动发展", "科技成果转化", 1),
    @("人工智能", "生物医药", 2),
    @("人工智能", "科技成果转化", 1),
    @("生物医药", "科技成果转化", 1),
    @("人工智能", "知识产权保护", 1),
    @("合成生物学", "生物医药", 1)
)

$r = 2
foreach ($row in $cooccurrenceData) {
    $wsCooccurrence.Cells.Item($r, 1).Value = $row[0]
    $wsCooccurrence.Cells.Item($r, 2).Value = $row[1]
    $wsCooccurrence.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Associations
# ---------------------------------------------------------------------------
$wsAssociations = $wb.Worksheets.Item("Associations")

$associationsData = @(
    @("新质生产力", 2),
    @("人工智能", 7),
    @("创新驱动发展", 1),
    @("科技成果转化", 2),
    @("量子通信", 1),
    @("生物医药", 4),
    @("国家创新体系", 1),
    @("合成生物学", 3),
    @("量子计算", 1),
    @("知识产权保护", 1),
    @("科技体制改革", 1)
)

$r = 2
foreach ($row in $associationsData) {
    $wsAssociations.Cells.Item($r, 1).Value = $row[0]
    $wsAssociations.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4) Header styling: match the bold/border/centered header style already
#    used by the "Associations" sheet (A1/B1), applying it to the
#    "Summary Table" header row and the new "Cooccurrence" header row.
# ---------------------------------------------------------------------------
$wsAssociations.Range("A1").Copy()
$wsSummary.Range("A1:E1").PasteSpecial(-4122)

$wsAssociations.Range("A1:B1").Copy()
$wsCooccurrence.Range("A1:C1").PasteSpecial(-4122)

$excel.CutCopyMode = 0
